$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet entry: 2020-07-28, 4 hours, long wrapped comment (new shared string).
$ws.Range("A41").Value = 44040
$ws.Range("B41").Value = 4
$ws.Range("C41").Value = "Nowa tabela bazie stworzona przy użyciu migracji EF. Zmieniona wiadomość wyjątku dodawania do bazy danych. `nLogowanie czasu dostarczenia raportu i czy został poprawnie dodany. Nowe testy dla XMLExport."

# Match formatting of similarly long wrapped comments elsewhere in column C.
$ws.Cells.Item(41, 3).WrapText = $true
$ws.Rows.Item(41).RowHeight = 31.5

# Move the selection to the next empty row, like the author did after typing the entry.
$ws.Range("C42").Select()
